$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: convert a cell that currently holds a *number* into one holding a
# textual label (e.g. "0" or "***.*"), reusing the formatting of a sibling
# cell ($styleRef) that is already a text-typed "N/A" cell so the resulting
# style index matches the rest of the sheet instead of minting a new one.
# ---------------------------------------------------------------------------
function Set-TextCell($range, $text, $styleRef) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $styleRef.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Helper: convert a cell that currently holds text back into a plain number,
# reusing the formatting of a sibling cell ($styleRef) that already carries
# the desired numeric style (counts vs. percentages).
# ---------------------------------------------------------------------------
function Set-NumCell($range, $number, $styleRef) {
    $styleRef.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $range.Value = $number
}

# ---------------------------------------------------------------------------
# Header text tweaks (stored as rich-text runs inside shared strings)
# ---------------------------------------------------------------------------

# "Volume 32   Number  10" -> "...  11"
$issueCell = $ws.Range("A8")
$issueFull = $issueCell.Value2
$issueStart = $issueFull.Length - 1
$issueChars = $issueCell.Characters($issueStart, 2)
$issueChars.Text = "11"

# "Report Covering the Week  3/3/2025  Through  3/9/2025"
#                          -> 3/10/2025            3/16/2025
$weekCell = $ws.Range("C9")
# Replace the right-hand date first so the left-hand date's fixed offset
# (27) is unaffected by the length change of the first replacement.
$weekChars2 = $weekCell.Characters(46, 8)
$weekChars2.Text = "3/16/2025"
$weekChars1 = $weekCell.Characters(27, 8)
$weekChars1.Text = "3/10/2025"

# ---------------------------------------------------------------------------
# Row 15 (Rape): numbers -> "N/A" style text
# ---------------------------------------------------------------------------
Set-TextCell $ws.Range("C15") "0"     $ws.Range("A15")
Set-TextCell $ws.Range("D15") "0"     $ws.Range("A15")
Set-TextCell $ws.Range("E15") "***.*" $ws.Range("A15")

# ---------------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 16
$ws.Range("J16").Value = 24
$ws.Range("K16").Value = -33.333333333333
$ws.Range("L16").Value = -30.434782608695
$ws.Range("M16").Value = -58.974358974359
$ws.Range("N16").Value = -91.623036649214

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -84.615384615384
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -62.068965517241
$ws.Range("I17").Value = 32
$ws.Range("J17").Value = 58
$ws.Range("K17").Value = -44.827586206896
$ws.Range("L17").Value = -28.888888888888
$ws.Range("M17").Value = 68.421052631578
$ws.Range("N17").Value = -20

# ---------------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 27
$ws.Range("J18").Value = 32
$ws.Range("K18").Value = -15.625
$ws.Range("L18").Value = -55.737704918032
$ws.Range("M18").Value = -60.294117647058
$ws.Range("N18").Value = -90.689655172413

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 66.666666666666
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = -26.086956521739
$ws.Range("I19").Value = 94
$ws.Range("J19").Value = 118
$ws.Range("K19").Value = -20.338983050847
$ws.Range("L19").Value = -28.244274809160
$ws.Range("M19").Value = 5.617977528089
$ws.Range("N19").Value = -33.333333333333

# ---------------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------------
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -25
$ws.Range("I20").Value = 30
$ws.Range("J20").Value = 36
$ws.Range("K20").Value = -16.666666666666
$ws.Range("L20").Value = 15.384615384615
$ws.Range("M20").Value = -16.666666666666
$ws.Range("N20").Value = -91.549295774647

# ---------------------------------------------------------------------------
# Row 21 (TOTAL, bold)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -39.285714285714
$ws.Range("F21").Value = 71
$ws.Range("G21").Value = 107
$ws.Range("H21").Value = -33.644859813084
$ws.Range("I21").Value = 201
$ws.Range("J21").Value = 271
$ws.Range("K21").Value = -25.830258302583
$ws.Range("L21").Value = -31.399317406143
$ws.Range("M21").Value = -22.093023255814
$ws.Range("N21").Value = -80.390243902439

# ---------------------------------------------------------------------------
# Row 22 (Transit): D/E numbers -> "N/A" style text
# ---------------------------------------------------------------------------
Set-TextCell $ws.Range("D22") "0"     $ws.Range("A22")
Set-TextCell $ws.Range("E22") "***.*" $ws.Range("A22")

# ---------------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = 2.777777777777
$ws.Range("F24").Value = 176
$ws.Range("G24").Value = 119
$ws.Range("H24").Value = 47.899159663865
$ws.Range("I24").Value = 371
$ws.Range("J24").Value = 363
$ws.Range("K24").Value = 2.203856749311
$ws.Range("L24").Value = 7.536231884057
$ws.Range("M24").Value = 101.630434782609

# ---------------------------------------------------------------------------
# Row 25 (Retail Theft)
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = -15
$ws.Range("F25").Value = 109
$ws.Range("G25").Value = 72
$ws.Range("H25").Value = 51.388888888888
$ws.Range("I25").Value = 241
$ws.Range("J25").Value = 207
$ws.Range("K25").Value = 16.425120772946
$ws.Range("L25").Value = 18.719211822660

# ---------------------------------------------------------------------------
# Row 26 (Misd. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 49
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = 11.363636363636
$ws.Range("I26").Value = 116
$ws.Range("J26").Value = 143
$ws.Range("K26").Value = -18.881118881118
$ws.Range("L26").Value = 8.411214953271
$ws.Range("M26").Value = 58.904109589041

# ---------------------------------------------------------------------------
# Row 27 (UCR Rape*): numbers -> "N/A" style text
# ---------------------------------------------------------------------------
Set-TextCell $ws.Range("C27") "0"     $ws.Range("A27")
Set-TextCell $ws.Range("D27") "0"     $ws.Range("A27")
Set-TextCell $ws.Range("E27") "***.*" $ws.Range("A27")

# ---------------------------------------------------------------------------
# Row 28 (Other Sex Crimes): D/E/G/H "N/A" style text -> numbers
# ---------------------------------------------------------------------------
Set-NumCell $ws.Range("D28") 1   $ws.Range("F28")
Set-NumCell $ws.Range("E28") -100 $ws.Range("K28")
$ws.Range("F28").Value = 2
Set-NumCell $ws.Range("G28") 1   $ws.Range("F28")
Set-NumCell $ws.Range("H28") 100 $ws.Range("K28")
$ws.Range("J28").Value = 9
$ws.Range("K28").Value = 11.111111111111
